$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Лист1")

# Update the economics figures (H4, H6, H7) that back the pie chart
$ws.Range("H4").Value = 508250
$ws.Range("H6").Value = 304950
$ws.Range("H7").Value = 39200

# Update the active cell selection to H6
$ws.Range("H6").Select()
